# Experiment 1 graphs added - refresh the latency trace data with new
# measurements and tidy up the selection / axis scaling on the chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Updated latency samples (columns B:NewReno, C:Reno, D:Tahoe, E:Vegas) ---
$data = @{
    2  = @(20.928355412799998, 20.928355412799998, 20.928355412799998, 20.885040916499999)
    3  = @(20.928226950300001, 20.928226950300001, 20.928226950300001, 20.919912711399999)
    4  = @(20.927335340599999, 20.927335340599999, 20.927335340599999, 20.919632160399999)
    5  = @(20.932727418700001, 20.932727418700001, 20.932727418700001, 20.9207659574)
    6  = @(20.948023016200001, 20.948023016200001, 20.948023016200001, 20.927419602400001)
    7  = @(20.957212852000001, 20.9598500225,      20.9512462623,      20.964693304699999)
    8  = @(21.014766589299999, 21.029710230300001, 21.0197427821,      21.046566669000001)
    9  = @(22.120052334,       21.8806972366,       22.1609262785,      22.8543193167)
    10 = @(21.894913874699999, 21.9828509602,       21.9651539119,      22.9140970222)
    11 = @(23.3084147157,      23.2762754098,       23.291156199700001, 22.254171428599999)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}

# --- Chart axis rescale: tighten the value-axis minimum now that the data
#     sits in a narrower band ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valAx = $chart.Axes(2)
$valAx.MinimumScale = 20.5

# --- Selection moved off the old range onto Q10 ---
$ws.Range("Q10").Select()
